$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare formatting (wrap text, matching bold header font where applicable) for the
# three newly appended rows before filling them in, mirroring the existing rows.
$ws.Range("A19:D21").WrapText = $true
$ws.Range("A19:D19").RowHeight = 34
$ws.Range("A20:D20").RowHeight = 17
$ws.Range("A21:D21").RowHeight = 34

# 1) New feature variable name typed into the freshly added row
$ws.Range("A19").Value = "total_cpd_bin"

# 2) New "origin" column header
$ws.Range("D1").Value = "origin"
$ws.Range("D1").WrapText = $true
$ws.Range("D1").Font.Bold = $true

# 3) Fill the new "origin" column values, in the order they were authored
$ws.Range("D2:D11").Value = "original"
$ws.Range("D12:D14").Value = "external"
$ws.Range("D15:D21").Value = "calculated"
$ws.Range("D2:D21").WrapText = $true

# 4) Finish filling out the rest of the new rows
$ws.Range("B19").Value = "total_cpd binned into intervals of 5cpd up to 80cpd (i.e., the final bin is 80 to infinity)"
$ws.Range("C19").Value = "numeric"

$ws.Range("A20").Value = "cpd_bin_label"
$ws.Range("B20").Value = "sequential label of total_cpd_bin"
$ws.Range("C20").Value = "ordinal"

$ws.Range("A21").Value = "prp_change"
$ws.Range("B21").Value = "ratio of total_cpd / baseline_cpd. If total_cpd = 0, prp_change = -1"
$ws.Range("C21").Value = "numeric"

# Update the view: select D21 as the active cell, matching the final selection
# state left by the author.
$excel.ActiveWindow.ScrollRow = 4
[void]$ws.Range("D21").Select()
